# timelog.xlsx - add Spring Boot (backend) project rows and rename the
# first entry so it refers to the Front end explicitly now that there is
# a Backend ("Bäkkäri") track too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Projektin käynnistys" (generic "project startup") -> "Frontin käynnistys"
# now that work has also started on the backend.
$ws.Range("E4").Value2 = "Frontin käynnistys"

# Tidy up the formatting on the existing last two rows so they match the
# rest of the date / description columns.
$ws.Range("A8:A9").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("E8:E9").NumberFormat = $ws.Range("E4").NumberFormat
$ws.Range("E8:E9").HorizontalAlignment = $ws.Range("E4").HorizontalAlignment

# New row 10: kicking off the Spring Boot backend project.
$ws.Range("A10").Value2 = 46034
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("B10").Value2 = 0.65625
$ws.Range("B10").NumberFormat = $ws.Range("B9").NumberFormat
$ws.Range("C10").Value2 = 0.6875
$ws.Range("C10").NumberFormat = $ws.Range("C9").NumberFormat
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("D10").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E10").Value2 = "Bäkkärin käynnistys"
$ws.Range("E10").NumberFormat = $ws.Range("E9").NumberFormat
$ws.Range("E10").HorizontalAlignment = $ws.Range("E9").HorizontalAlignment

# New row 11: continuing backend programming (points request validation).
$ws.Range("A11").Value2 = 46035
$ws.Range("A11").NumberFormat = $ws.Range("A9").NumberFormat
$ws.Range("B11").Value2 = 0.770833333333333
$ws.Range("B11").NumberFormat = $ws.Range("B9").NumberFormat
$ws.Range("C11").Value2 = 0.847222222222222
$ws.Range("C11").NumberFormat = $ws.Range("C9").NumberFormat
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("D11").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E11").Value2 = "Bäkkäri ohjelmointi"
$ws.Range("E11").NumberFormat = $ws.Range("E3").NumberFormat
$ws.Range("E11").HorizontalAlignment = $ws.Range("E3").HorizontalAlignment

# Match the author's final cursor position.
$ws.Range("C11").Select()
